# Restore C10 (From value for rule R30) back to 1, matching the
# "Restored from revision" commit that reverted this cell's value.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$ws.Range("C10").Value = 1
